$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9074
$ws.Range("I18").Value = 9829.5
$ws.Range("J18").Value = 4541
$ws.Range("K18").Value = 9829.5
$ws.Range("L18").Value = 4541
$ws.Range("M18").Value = -9545.5
$ws.Range("N18").Value = -5109

$ws.Range("H28").Value = 831.6087
$ws.Range("I28").Value = 923.7778
$ws.Range("K28").Value = 923.7778
$ws.Range("M28").Value = -438.7778

$ws.Range("H41").Value = 846.6316
$ws.Range("I41").Value = 1010.7143
$ws.Range("K41").Value = 1010.7143
$ws.Range("M41").Value = -570.7143

$ws.Range("H46").Value = 1001149.4
$ws.Range("I46").Value = 497.5
$ws.Range("K46").Value = 1492.5
$ws.Range("M46").Value = -1373.5

$ws.Range("H60").Value = 1001149.4
$ws.Range("I60").Value = 497.5
$ws.Range("K60").Value = 1492.5
$ws.Range("M60").Value = -1008.5

$ws.Range("H99").Value = 125021070
$ws.Range("I99").Value = 27932
$ws.Range("J99").Value = 500000500
$ws.Range("K99").Value = 83796
$ws.Range("L99").Value = 1500001500
$ws.Range("M99").Value = -82298
$ws.Range("N99").Value = -1500004496

$ws.Range("H107").Value = 40214.21
$ws.Range("I107").Value = 592.6667
$ws.Range("K107").Value = 592.6667
$ws.Range("M107").Value = 1327.3333

$ws.Range("H138").Value = 2424.932
$ws.Range("I138").Value = 1427.5333
$ws.Range("K138").Value = 4282.5999
$ws.Range("M138").Value = 857.4000999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4925.7925
$ws.Range("I32").Value = 4181.8335
$ws.Range("K32").Value = 4181.8335
$ws.Range("M32").Value = -3894.8335

$ws.Range("H45").Value = 7027.44
$ws.Range("I45").Value = 9777.538
$ws.Range("J45").Value = 4048.1667
$ws.Range("K45").Value = 9777.538
$ws.Range("L45").Value = 4048.1667
$ws.Range("M45").Value = -9400.538
$ws.Range("N45").Value = -4802.1667

$ws.Range("H61").Value = 2090.3157
$ws.Range("I61").Value = 1937.8125
$ws.Range("J61").Value = 2903.6667
$ws.Range("K61").Value = 1937.8125
$ws.Range("L61").Value = 2903.6667
$ws.Range("M61").Value = -1725.8125
$ws.Range("N61").Value = -3327.6667

$ws.Range("H74").Value = 2085.25
$ws.Range("I74").Value = 1808.6154
$ws.Range("K74").Value = 1808.6154
$ws.Range("M74").Value = -934.6153999999999

$ws.Range("H77").Value = 2085.25
$ws.Range("I77").Value = 1808.6154
$ws.Range("K77").Value = 9043.076999999999
$ws.Range("M77").Value = -4675.076999999999

$ws.Range("H110").Value = 2332.1428
$ws.Range("I110").Value = 1440.6875
$ws.Range("K110").Value = 1440.6875
$ws.Range("M110").Value = 604.3125

$ws.Range("H136").Value = 2090.3157
$ws.Range("I136").Value = 1937.8125
$ws.Range("J136").Value = 2903.6667
$ws.Range("K136").Value = 5813.4375
$ws.Range("L136").Value = 8711.000100000001
$ws.Range("M136").Value = -3263.4375
$ws.Range("N136").Value = -13811.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5551.25
$ws.Range("I20").Value = 7088.778
$ws.Range("J20").Value = 2783.7
$ws.Range("K20").Value = 7088.778
$ws.Range("L20").Value = 2783.7
$ws.Range("M20").Value = -6841.778
$ws.Range("N20").Value = -3277.7

$ws.Range("H86").Value = 3041.3635
$ws.Range("I86").Value = 3059.4119
$ws.Range("K86").Value = 3059.4119
$ws.Range("M86").Value = -1936.4119

$ws.Range("H89").Value = 3041.3635
$ws.Range("I89").Value = 3059.4119
$ws.Range("K89").Value = 15297.0595
$ws.Range("M89").Value = -9681.059499999999

$ws.Range("H94").Value = 774.4706
$ws.Range("I94").Value = 791.38464
$ws.Range("J94").Value = 719.5
$ws.Range("K94").Value = 791.38464
$ws.Range("L94").Value = 719.5
$ws.Range("M94").Value = -340.38464
$ws.Range("N94").Value = -1621.5

$ws.Range("H99").Value = 9084.037
$ws.Range("I99").Value = 10455.392
$ws.Range("J99").Value = 1198.75
$ws.Range("K99").Value = 10455.392
$ws.Range("L99").Value = 1198.75
$ws.Range("M99").Value = -8957.392
$ws.Range("N99").Value = -4194.75

$ws.Range("H134").Value = 1617.3914
$ws.Range("I134").Value = 1252.6842
$ws.Range("K134").Value = 3758.0526
$ws.Range("M134").Value = -1223.0526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4831.1577
$ws.Range("I31").Value = 3310.3
$ws.Range("K31").Value = 3310.3
$ws.Range("M31").Value = -3015.3

$ws.Range("H34").Value = 4831.1577
$ws.Range("I34").Value = 3310.3
$ws.Range("K34").Value = 3310.3
$ws.Range("M34").Value = -3108.3

$ws.Range("H52").Value = 101250.664
$ws.Range("J52").Value = 124387.5
$ws.Range("L52").Value = 124387.5
$ws.Range("N52").Value = -124975.5

$ws.Range("H69").Value = 9000
$ws.Range("I69").Value = 9000
$ws.Range("K69").Value = 9000
$ws.Range("M69").Value = -8251

$ws.Range("H72").Value = 9000
$ws.Range("I72").Value = 9000
$ws.Range("K72").Value = 27000
$ws.Range("M72").Value = -23256

$ws.Range("H132").Value = 4189.1304
$ws.Range("I132").Value = 4978.25
$ws.Range("K132").Value = 14934.75
$ws.Range("M132").Value = -12404.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2198.6
$ws.Range("I63").Value = 2198.6
$ws.Range("K63").Value = 6595.799999999999
$ws.Range("M63").Value = -5846.799999999999

$ws.Range("H66").Value = 2198.6
$ws.Range("I66").Value = 2198.6
$ws.Range("K66").Value = 19787.4
$ws.Range("M66").Value = -16043.4

$ws.Range("H70").Value = 103147.1
$ws.Range("I70").Value = 144568.42
$ws.Range("K70").Value = 433705.26
$ws.Range("M70").Value = -433390.26

$ws.Range("H73").Value = 103147.1
$ws.Range("I73").Value = 144568.42
$ws.Range("K73").Value = 433705.26
$ws.Range("M73").Value = -432613.26

$ws.Range("H88").Value = 12499
$ws.Range("J88").Value = 12499
$ws.Range("L88").Value = 37497
$ws.Range("N88").Value = -38353

$ws.Range("H91").Value = 12499
$ws.Range("J91").Value = 12499
$ws.Range("L91").Value = 37497
$ws.Range("N91").Value = -40461

$ws.Range("H94").Value = 7323.1665
$ws.Range("I94").Value = 5646.3335
$ws.Range("K94").Value = 16939.0005
$ws.Range("M94").Value = -16263.0005

$ws.Range("H106").Value = 5765.3184
$ws.Range("J106").Value = 5941.95
$ws.Range("L106").Value = 17825.85
$ws.Range("N106").Value = -19717.85

$ws.Range("H114").Value = 25000974
$ws.Range("I114").Value = 28572256
$ws.Range("K114").Value = 85716768
$ws.Range("M114").Value = -85713514

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 89319.75
$ws.Range("J32").Value = 89093
$ws.Range("L32").Value = 89093
$ws.Range("N32").Value = -89685

$ws.Range("H33").Value = 24999.666

$ws.Range("H97").Value = 27794.861
$ws.Range("I97").Value = 39322.75
$ws.Range("J97").Value = 2177.3333
$ws.Range("K97").Value = 39322.75
$ws.Range("L97").Value = 2177.3333
$ws.Range("M97").Value = -38826.75
$ws.Range("N97").Value = -3169.3333

$ws.Range("H132").Value = 4422.8696
$ws.Range("I132").Value = 5284.5713
$ws.Range("K132").Value = 15853.7139
$ws.Range("M132").Value = -13323.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 799.36365
$ws.Range("I16").Value = 808.75
$ws.Range("K16").Value = 808.75
$ws.Range("M16").Value = -638.75

$ws.Range("H46").Value = 22523.955
$ws.Range("I46").Value = 35675.46
$ws.Range("J46").Value = 3527.3333
$ws.Range("K46").Value = 35675.46
$ws.Range("L46").Value = 3527.3333
$ws.Range("M46").Value = -35487.46
$ws.Range("N46").Value = -3903.3333

$ws.Range("H55").Value = 54156
$ws.Range("I55").Value = 69667.47
$ws.Range("J55").Value = 7621.6
$ws.Range("K55").Value = 69667.47
$ws.Range("L55").Value = 7621.6
$ws.Range("M55").Value = -69494.47
$ws.Range("N55").Value = -7967.6

$ws.Range("H68").Value = 1550
$ws.Range("I68").Value = 1437.5
$ws.Range("K68").Value = 1437.5
$ws.Range("M68").Value = -688.5

$ws.Range("H71").Value = 1550
$ws.Range("I71").Value = 1437.5
$ws.Range("K71").Value = 7187.5
$ws.Range("M71").Value = -3443.5

$ws.Range("H93").Value = 22061.37
$ws.Range("I93").Value = 4453.091
$ws.Range("K93").Value = 4453.091
$ws.Range("M93").Value = -3205.091

$ws.Range("H132").Value = 4409.5454
$ws.Range("I132").Value = 4381.4546
$ws.Range("J132").Value = 4465.727
$ws.Range("K132").Value = 13144.3638
$ws.Range("L132").Value = 13397.181
$ws.Range("M132").Value = -10614.3638
$ws.Range("N132").Value = -18457.181

$ws.Range("H136").Value = 3353.7778
$ws.Range("I136").Value = 3198.0715
$ws.Range("K136").Value = 9594.2145
$ws.Range("M136").Value = -7044.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 27691.584
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H46").Value = 52999.332
$ws.Range("J46").Value = 52999.332
$ws.Range("L46").Value = 52999.332
$ws.Range("N46").Value = -53461.332

$ws.Range("H134").Value = 52999.332
$ws.Range("J134").Value = 52999.332
$ws.Range("L134").Value = 158997.996
$ws.Range("N134").Value = -164067.996
